$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()
$ws.Range("G1").Value = "step_values"
$ws.Range("G1").Select()
